$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 29.17403400000001
$ws.Range("N2").Value = 87.52210200000002
$ws.Range("O2").Value = 0.3835306213274714
$ws.Range("P2").Value = 0.3835306213274714
$ws.Range("Q2").Value = 5.183175576576001
$ws.Range("R2").Value = 46.64858018918401
$ws.Range("S2").Value = 0.3835306213274714
$ws.Range("T2").Value = 0.3835306213274714

# Row 3 updates
$ws.Range("O3").Value = 0.5274816184042599
$ws.Range("P3").Value = 0.5274816184042599
$ws.Range("S3").Value = 0.5274816184042599
$ws.Range("T3").Value = 0.5274816184042599

# Row 4 updates
$ws.Range("M4").Value = 6.769034333333334
$ws.Range("N4").Value = 20.307103
$ws.Range("O4").Value = 0.08898776026826866
$ws.Range("P4").Value = 0.08898776026826867
$ws.Range("Q4").Value = 1.202613715797334
$ws.Range("R4").Value = 10.823523442176
$ws.Range("S4").Value = 0.08898776026826866
$ws.Range("T4").Value = 0.08898776026826867
